# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de), the a.md / b.md rows are now
# "handed back" and in sync with en-US:
#   - Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns
#     get populated (previously empty) with hyperlinks to the same md /
#     xlf files referenced by the handoff columns (A / C)
#   - The "Latest Handback DateTime" (G) is stamped with the handback time
#
# The ".localization-config" row (row 4) is "Not to be localized" / not
# handed off, so it is left untouched.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$localeInfo = @{
    "zh-cn" = @{
        XlfName        = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandoffXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e539883cfd49214c05963e347f9c26d185e9b2ca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandbackTime   = "2016-02-23 07:29:10"
    }
    "de-de" = @{
        XlfName        = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandoffXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baeb1a34094d07e71a8ac46d838a16fd45085b00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandbackTime   = "2016-02-23 07:29:35"
    }
}

$mdUrls = @{
    "a.md" = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/a.md"
    "b.md" = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/b.md"
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $localeInfo[$sheetName]

    # Rows 2 (a.md) and 3 (b.md) both get handed back.
    foreach ($row in @(2, 3)) {
        $sourceName = $ws.Range("A" + $row).Text

        # B: Status -> handed back, in sync with en-US
        $ws.Range("B" + $row).Value = $newStatus

        # E: Latest Target File -> same source file (a.md), now hyperlinked
        $ws.Hyperlinks.Add($ws.Range("E" + $row), $mdUrls["a.md"], "", "", "a.md") | Out-Null

        # F: Latest Handback File -> the handoff xlf, now hyperlinked
        $ws.Hyperlinks.Add($ws.Range("F" + $row), $info.HandoffXlfUrl, "", "", $info.XlfName) | Out-Null

        # G: Latest Handback DateTime -> stamp of this handback run
        $ws.Range("G" + $row).Value = $info.HandbackTime
    }
}
